# Manuscript III analysis update - "19812000" worksheet
# Adds median/IQR columns (replacing average/Std headers), relabels the
# treatment-group rows with their text names, adds a KS D-stat / Cliffs D
# comparison block (row 5) and fills in the corresponding KS Test row (row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("19812000")

# --- Row 1 headers: average/Std -> median/IQR --------------------------
$ws.Range("B1").Value = "Duration (s), median"
$ws.Range("C1").Value = "Duration (s), IQR"

# --- Row 2: "Control" treatment group -----------------------------------
$ws.Range("A2").Value = "Control"
$ws.Range("B2").Value = 24.681840000000136
$ws.Range("C2").Value = 5.4134399999998237
$ws.Range("D2").Value = 0.2243857360774992
$ws.Range("E2").Value = 391781.4454027852
$ws.Range("F2").Value = 346234.2190919249

# --- Row 3: "Test" treatment group ---------------------------------------
$ws.Range("A3").Value = "Test"
$ws.Range("B3").Value = 26.859647999999652
$ws.Range("C3").Value = 2.9741279999996095
$ws.Range("E3").Value = 426501.76747243921
$ws.Range("F3").Value = 128338.2572215955

# --- Row 4: "Washout" treatment group -------------------------------------
$ws.Range("A4").Value = "Washout"
$ws.Range("B4").Value = 24.66172799999913
$ws.Range("C4").Value = 5.5079279999995379
$ws.Range("E4").Value = 304965.69681405724
$ws.Range("F4").Value = 73571.530899229925

# --- Row 5 (new): KS D stat / Cliffs D column headers ---------------------
$ws.Range("B5").Value = "p-value"
$ws.Range("C5").Value = "KS D stat"
$ws.Range("D5").Value = "Cliffs D"
$ws.Range("E5").Value = "p-value"
$ws.Range("F5").Value = "KS D stat"
$ws.Range("G5").Value = "Cliffs D"

# --- Row 6: KS Test, 1 vs 2 results, with new KS D stat/Cliffs D values ---
$ws.Range("C6").Value = 0.40367965367965369
$ws.Range("D6").Value = -0.44805194805194803
$ws.Range("F6").Value = 0.36255411255411252
$ws.Range("G6").Value = -0.084415584415584416

# --- View state: scrolled down to the bottom comparison table -------------
$ws.Range("A28:F30").Select() | Out-Null
